$wb = $excel.ActiveWorkbook

$changes = @{
    5  = 24
    6  = 116
    8  = 464
    9  = 50
    11 = 581
    13 = 305
    20 = 51
    21 = 101
    22 = 960
    23 = 1407
    24 = 302
    25 = 332
    27 = 77
    30 = 87
    32 = 256
    33 = 276
    34 = 1627
    35 = 53
    38 = 587
    40 = 3710
    41 = 432
    42 = 205
    43 = 915
    46 = 72
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $changes.Keys) {
        $ws.Cells.Item($row, 6).Value = $changes[$row]
    }
}
